$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 2 (Oct 25 2020 match) so that the former row 3
# (Oct 10 2020 match) shifts up and becomes the new row 2.
$ws.Rows(2).Delete()
